$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 12: Polygon -> WrappedEther (updated data)
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "1.881.49"
$ws.Cells.Item(12, 5).Value = "  -0.61%  "

# Row 13: WrappedEther -> Polygon (updated data)
$ws.Cells.Item(13, 2).Value = "Polygon"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Cells.Item(13, 4) "0.7240"
$ws.Cells.Item(13, 5).Value = "  -2.17%  "

# Update Price (D) / Volume(1h) (E) for the remaining rows
$ws.Cells.Item(2, 4).Value = "29.104.89"
$ws.Cells.Item(2, 5).Value = "  -2.19%  "
$ws.Cells.Item(3, 4).Value = "1.848.72"
$ws.Cells.Item(3, 5).Value = "  -1.05%  "
$ws.Cells.Item(4, 5).Value = "  +0.21%  "
Set-TextValue $ws.Cells.Item(5, 4) "0.6923"
$ws.Cells.Item(5, 5).Value = "  -5.16%  "
Set-TextValue $ws.Cells.Item(6, 4) "237.82"
$ws.Cells.Item(6, 5).Value = "  -1.21%  "
Set-TextValue $ws.Cells.Item(7, 4) "1.001"
$ws.Cells.Item(7, 5).Value = "  +0.13%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.07704"
$ws.Cells.Item(8, 5).Value = "  +8.60%  "
$ws.Cells.Item(9, 5).Value = "  -2.93%  "
Set-TextValue $ws.Cells.Item(10, 4) "23.20"
$ws.Cells.Item(10, 5).Value = "  -4.45%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.08112"
$ws.Cells.Item(11, 5).Value = "  -1.13%  "
Set-TextValue $ws.Cells.Item(14, 4) "5.206"
$ws.Cells.Item(14, 5).Value = "  -2.11%  "
Set-TextValue $ws.Cells.Item(15, 4) "88.91"
$ws.Cells.Item(15, 5).Value = "  -3.35%  "
$ws.Cells.Item(16, 4).Value = "29.107.91"
$ws.Cells.Item(16, 5).Value = "  -2.19%  "
$ws.Cells.Item(17, 5).Value = "  -4.42%  "
Set-TextValue $ws.Cells.Item(18, 4) "0.000007774"
$ws.Cells.Item(18, 5).Value = "  -0.21%  "
$ws.Cells.Item(19, 5).Value = "  -1.21%  "
Set-TextValue $ws.Cells.Item(20, 4) "235.27"
$ws.Cells.Item(20, 5).Value = "  -4.95%  "
Set-TextValue $ws.Cells.Item(21, 4) "1.000"
$ws.Cells.Item(21, 5).Value = "  +0.10%  "
$ws.Cells.Item(22, 4).Value = "2.095.48"
$ws.Cells.Item(22, 5).Value = "  -1.55%  "
Set-TextValue $ws.Cells.Item(23, 4) "1.002"
$ws.Cells.Item(23, 5).Value = "  +0.27%  "
Set-TextValue $ws.Cells.Item(24, 4) "7.594"
$ws.Cells.Item(24, 5).Value = "  -1.80%  "
Set-TextValue $ws.Cells.Item(25, 4) "8.966"
$ws.Cells.Item(25, 5).Value = "  -2.26%  "
Set-TextValue $ws.Cells.Item(26, 4) "160.97"
$ws.Cells.Item(26, 5).Value = "  -1.33%  "
Set-TextValue $ws.Cells.Item(27, 4) "0.1431"
$ws.Cells.Item(27, 5).Value = "  -7.15%  "
Set-TextValue $ws.Cells.Item(28, 4) "18.00"
$ws.Cells.Item(28, 5).Value = "  -2.56%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.978"
$ws.Cells.Item(29, 5).Value = "  -1.26%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.399"
$ws.Cells.Item(30, 5).Value = "  -2.78%  "
Set-TextValue $ws.Cells.Item(31, 4) "4.484"
$ws.Cells.Item(31, 5).Value = "  -0.65%  "
$ws.Cells.Item(32, 5).Value = "  -2.30%  "
Set-TextValue $ws.Cells.Item(33, 4) "4.015"
$ws.Cells.Item(33, 5).Value = "  -3.52%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.05218"
$ws.Cells.Item(34, 5).Value = "  -1.13%  "
$ws.Cells.Item(35, 5).Value = "  -3.91%  "
$ws.Cells.Item(36, 5).Value = "  +2.17%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.6984"
$ws.Cells.Item(37, 5).Value = "  -6.00%  "
Set-TextValue $ws.Cells.Item(38, 4) "2.658"
$ws.Cells.Item(38, 5).Value = "  -1.16%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.01848"
$ws.Cells.Item(39, 5).Value = "  -4.11%  "
Set-TextValue $ws.Cells.Item(40, 4) "2.680"
$ws.Cells.Item(40, 5).Value = "  -1.85%  "
Set-TextValue $ws.Cells.Item(41, 4) "0.9147"
$ws.Cells.Item(41, 5).Value = "  +5.43%  "
Set-TextValue $ws.Cells.Item(42, 4) "5.985"
$ws.Cells.Item(42, 5).Value = "  +0.06%  "
$ws.Cells.Item(43, 4).Value = "1.078.83"
$ws.Cells.Item(43, 5).Value = "  +3.48%  "
Set-TextValue $ws.Cells.Item(44, 4) "0.4251"
$ws.Cells.Item(44, 5).Value = "  -4.33%  "
Set-TextValue $ws.Cells.Item(45, 4) "70.45"
$ws.Cells.Item(45, 5).Value = "  -0.99%  "
Set-TextValue $ws.Cells.Item(46, 4) "1.001"
$ws.Cells.Item(46, 5).Value = "  +0.13%  "
Set-TextValue $ws.Cells.Item(47, 4) "103.29"
$ws.Cells.Item(47, 5).Value = "  -0.37%  "
Set-TextValue $ws.Cells.Item(48, 4) "1.766"
$ws.Cells.Item(48, 5).Value = "  -2.40%  "
$ws.Cells.Item(49, 4).Value = "1.991.78"
$ws.Cells.Item(49, 5).Value = "  -1.71%  "
Set-TextValue $ws.Cells.Item(50, 4) "9.115"
$ws.Cells.Item(50, 5).Value = "  -4.01%  "
Set-TextValue $ws.Cells.Item(51, 4) "6.987"
